$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")
$values = @(3125.119287316165, 0, 38683.01019569611, 0, 231779.209144148, 7546.507577201783, 0, 2027.422343033684, 0, 0, 0, 0, 0, 1894.610970861271, 1596.959756313929)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
    $ws.Cells.Item(2, $i + 1).Value2 = $values[$i]
}

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")
$values = @(5592.841625345535, 0, 158331.0001646095, 0, 231779.209144148, 13162.01156291295, 0, 6555.616545050024, 0, 0, 0, 0, 0, 6037.449078019277, 5007.288439540145)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
    $ws.Cells.Item(2, $i + 1).Value2 = $values[$i]
}

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")
$values = @(24989.03564310195, 0, 233797.8217711074, 0, 231779.209144148, 13276.08564128262, 0, 9705.535366320997, 0, 0, 0, 0, 0, 10316.13445312728, 7412.70627513708)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
    $ws.Cells.Item(2, $i + 1).Value2 = $values[$i]
}

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")
$values = @(24989.03564310195, 0, 233797.8217711074, 0, 231779.209144148, 13276.08564128262, 0, 9705.535366320997, 0, 0, 0, 0, 0, 11242.41425508182, 7412.70627513708)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
    $ws.Cells.Item(2, $i + 1).Value2 = $values[$i]
}

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")
$values = @(31125.48275843248, 154.4757918726473, 233797.8217711074, 0, 231779.209144148, 13276.08564128262, 0, 9705.535366320997, 0, 0, 0, 0, 0, 13513.47261752928, 8079.009742155169)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
    $ws.Cells.Item(2, $i + 1).Value2 = $values[$i]
}

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")
$values = @(31125.48275843248, 154.4757918726473, 233797.8217711074, 0, 231779.209144148, 13276.08564128262, 0, 9705.535366320997, 0, 0, 0, 0, 0, 13513.47261752928, 8079.009742155169)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
    $ws.Cells.Item(2, $i + 1).Value2 = $values[$i]
}
